$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.986.23'
$ws.Range('E2').Value = '  +3.40%  '
$ws.Range('D3').Value = '1.726.46'
$ws.Range('E3').Value = '  +3.07%  '
$ws.Range('D4').Value = "'1.00"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = "'219.14"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.81%  '
$ws.Range('E6').Value = '  +1.53%  '
$ws.Range('D7').Value = "'1.00"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').Value = "'24.17"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +13.85%  '
$ws.Range('E9').Value = '  +3.76%  '
$ws.Range('E10').Value = '  +2.33%  '
$ws.Range('D11').Value = "'0.0902"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.10%  '
$ws.Range('D12').Value = '1.971.90'
$ws.Range('E12').Value = '  +3.17%  '
$ws.Range('D13').Value = '1.717.82'
$ws.Range('E13').Value = '  +2.56%  '
$ws.Range('D14').Value = "'4.26"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.57%  '
$ws.Range('D15').Value = "'0.566"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +6.13%  '
$ws.Range('D16').Value = "'67.89"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.86%  '
$ws.Range('D17').Value = '27.935.69'
$ws.Range('E17').Value = '  +3.30%  '
$ws.Range('D18').Value = "'243.26"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.47%  '
$ws.Range('E19').Value = '  +2.47%  '
$ws.Range('D20').Value = "'7.91"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.81%  '
$ws.Range('D21').Value = "'0.999"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.14%  '
$ws.Range('E22').Value = '  +4.59%  '
$ws.Range('E23').Value = '  +4.82%  '
$ws.Range('E24').Value = '  +0.12%  '
$ws.Range('D25').Value = "'149.13"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.03%  '
$ws.Range('D26').Value = "'7.53"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +4.22%  '
$ws.Range('D27').Value = "'16.82"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.79%  '
$ws.Range('E29').Value = '  -0.15%  '
$ws.Range('D30').Value = "'0.0512"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.81%  '
$ws.Range('E31').Value = '  +1.89%  '
$ws.Range('D32').Value = "'3.45"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.01%  '
$ws.Range('E33').Value = '  +3.07%  '
$ws.Range('D34').Value = '1.486.02'
$ws.Range('E34').Value = '  -3.83%  '
$ws.Range('E35').Value = '  -1.69%  '
$ws.Range('D36').Value = "'0.961"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +4.25%  '
$ws.Range('D37').Value = "'0.611"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.02%  '
$ws.Range('D39').Value = "'0.0175"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.93%  '
$ws.Range('E40').Value = '  +1.10%  '
$ws.Range('D41').Value = "'71.43"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +5.60%  '
$ws.Range('E42').Value = '  +4.52%  '
$ws.Range('D43').Value = "'1.00"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.08%  '
$ws.Range('E44').Value = '  +1.72%  '
$ws.Range('D45').Value = '1.874.91'
$ws.Range('E45').Value = '  +3.11%  '
$ws.Range('E46').Value = '  +1.19%  '
$ws.Range('E47').Value = '  +13.84%  '
$ws.Range('D48').Value = "'91.88"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.08%  '
$ws.Range('E49').Value = '  +4.05%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = "'8.29"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.03%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').Value = "'0.106"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.33%  '
